# Generate Report for Handoff
# Adds two new localization files (479097ab-... and e99aeb79-...) to the
# Overview / zh-cn / de-de sheets, each contributing one new row.

$wb = $excel.ActiveWorkbook

# ---- constants describing the two new files -----------------------------
$guid1 = "479097ab-a4e8-4268-9123-34a164937545"
$guid2 = "e99aeb79-1668-471b-a8cd-51eb14b9cafd"

$zhHash1 = "0565c5870c713f7020e5bf74b6dee2bfe2ffe3f9"
$zhHash2 = "50bb61a3458907cecf0cd1d7c5ae9246025e9f1b"

$status = "Ready for handoff"
$overviewDate = "2016-36-21 04:36:25"
$handoffDate = "2016-03-21 04:36:21"
$handoffDateDe = "2016-03-21 04:36:25"
$noHandback = "0001-01-01 00:00:00"
$reason = "Include"

$file1md = "$guid1.md"
$file2md = "$guid2.md"

$file1zh = "$guid1.$zhHash1.zh-cn.xlf"
$file2zh = "$guid2.$zhHash2.zh-cn.xlf"
$file1de = "$guid1.$zhHash1.de-de.xlf"
$file2de = "$guid2.$zhHash2.de-de.xlf"

# Hyperlink-style font: same underline/colour as the existing links so the
# style dedupes against the workbook's pre-existing "HyperLink" font.
$hlUnderline = 2          # xlUnderlineStyleSingle
$hlColor = 15570276       # 0x00ED9564 == RGB(100,149,237) == #6495ED

function Apply-LinkFont($range) {
    $f = $range.Font
    $f.Underline = $hlUnderline
    $f.Color = $hlColor
}

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value2 = $file1md
$ws1.Range("B4").Value2 = $status
$ws1.Range("C4").Value2 = $status
$ws1.Range("D4").Value2 = $overviewDate

$ws1.Range("A5").Value2 = $file2md
$ws1.Range("B5").Value2 = $status
$ws1.Range("C5").Value2 = $status
$ws1.Range("D5").Value2 = $overviewDate

$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file1md", "", "", $file1md)
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file2md", "", "", $file2md)

Apply-LinkFont $ws1.Range("A4")
Apply-LinkFont $ws1.Range("A5")

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value2 = $file1md
$ws2.Range("B4").Value2 = ".md"
$ws2.Range("C4").Value2 = $status
$ws2.Range("D4").Value2 = $file1zh
$ws2.Range("E4").Value2 = $handoffDate
$ws2.Range("H4").Value2 = $noHandback
$ws2.Range("I4").Value2 = $reason

$ws2.Range("A5").Value2 = $file2md
$ws2.Range("B5").Value2 = ".md"
$ws2.Range("C5").Value2 = $status
$ws2.Range("D5").Value2 = $file2zh
$ws2.Range("E5").Value2 = $handoffDate
$ws2.Range("H5").Value2 = $noHandback
$ws2.Range("I5").Value2 = $reason

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file1md", "", "", $file1md)
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$file1md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file1zh", "", "", $file1zh)

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file2md", "", "", $file2md)
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/master/e2e/$file2md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$file2zh", "", "", $file2zh)

Apply-LinkFont $ws2.Range("A4")
Apply-LinkFont $ws2.Range("B4")
Apply-LinkFont $ws2.Range("D4")
Apply-LinkFont $ws2.Range("A5")
Apply-LinkFont $ws2.Range("B5")
Apply-LinkFont $ws2.Range("D5")

$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value2 = $file1md
$ws3.Range("B4").Value2 = ".md"
$ws3.Range("C4").Value2 = $status
$ws3.Range("D4").Value2 = $file1de
$ws3.Range("E4").Value2 = $handoffDateDe
$ws3.Range("H4").Value2 = $noHandback
$ws3.Range("I4").Value2 = $reason

$ws3.Range("A5").Value2 = $file2md
$ws3.Range("B5").Value2 = ".md"
$ws3.Range("C5").Value2 = $status
$ws3.Range("D5").Value2 = $file2de
$ws3.Range("E5").Value2 = $handoffDateDe
$ws3.Range("H5").Value2 = $noHandback
$ws3.Range("I5").Value2 = $reason

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file1md", "", "", $file1md)
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$file1md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file1de", "", "", $file1de)

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file2md", "", "", $file2md)
$ws3.Hyperlinks.Add($ws3.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/master/e2e/$file2md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$file2de", "", "", $file2de)

Apply-LinkFont $ws3.Range("A4")
Apply-LinkFont $ws3.Range("B4")
Apply-LinkFont $ws3.Range("D4")
Apply-LinkFont $ws3.Range("A5")
Apply-LinkFont $ws3.Range("B5")
Apply-LinkFont $ws3.Range("D5")

$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Handback rows added for $guid1 and $guid2"
